# CRMAccuracyData.xlsx - append new CRM accuracy readings taken 2021-03-14
# (CRM opened 20210228, Batch # 141) and refresh the sheet's scroll/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Date, CRM value, Batch value triples for rows 12-19
$newRows = @(
    @(20210314, 2182.6919702618002, 2235.0700000000002),
    @(20210314, 2183.0707049272701, 2236.0700000000002),
    @(20210314, 2179.9893168764502, 2237.0700000000002),
    @(20210314, 2180.1189155137999, 2238.0700000000002),
    @(20210314, 2194.2628893932401, 2239.0700000000002),
    @(20210314, 2194.8087986166001, 2240.0700000000002),
    @(20210314, 2192.1340883392099, 2241.0700000000002),
    @(20210314, 2192.25361055657,   2242.0700000000002)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]                 # A: Date
    $ws.Cells.Item($r, 2).Value = $vals[1]                 # B: CRM value
    $ws.Cells.Item($r, 3).Value = $vals[2]                 # C: Batch value
    $ws.Cells.Item($r, 4).Formula = "=100*(B$r-C$r)/C$r"   # D: % off
    $ws.Cells.Item($r, 5).Value = 141                      # E: Batch #
    $ws.Cells.Item($r, 6).Value = "CRM opened 20210228"    # F: Notes
}

# Scroll the view down so row 10 is at the top and refresh the active selection
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("H15").Select()
